# Update the daily COVID-19 Valais figures for rows 269-281 (dates
# 27.11.2020 through 09.12.2020, serials A44155..A44167) on worksheet
# "Feuil1".
#
# Columns B (cumul cas positifs), H (total hospitalisations COVID-19),
# J (cumul decces COVID-19) and K (nb nouveaux decces COVID-19) are live
# formulas (e.g. "=IF(TODAY()>A268,B268+C269,"")") that recompute
# automatically once the underlying input cells are written, so only the
# input columns - C (nb nouveaux cas positifs), E (patients SI),
# F (intubes), G (hospit. hors SI), L (decces hopital) et M (decces
# extra-hospitaliers) - need to be set here; row 282's formulas then also
# recompute on their own once row 281 is filled in.
#
# Columns L and M are formatted as Text ("@") cells, so a plain .Value
# write would store the new figure as a text string instead of a number
# (mirroring how Excel treats a Text-formatted cell). Flip the cell to
# General, write the number, then restore the original Text format -
# exactly what happens when Excel itself enters a number that still
# displays/validates fine in a "@" cell but is stored as a real number.

function Set-NumericValue($cell, $value) {
    $originalFormat = $cell.NumberFormat
    $cell.NumberFormat = "General"
    $cell.Value = $value
    $cell.NumberFormat = $originalFormat
}

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Row 269 (27.11.2020)
$ws.Cells.Item(269, 3).Value = 152                      # C269
Set-NumericValue $ws.Cells.Item(269, 13) 4               # M269

# Row 272 (30.11.2020)
Set-NumericValue $ws.Cells.Item(272, 13) 2                # M272

# Row 273 (01.12.2020)
Set-NumericValue $ws.Cells.Item(273, 13) 9                # M273

# Row 275 (03.12.2020)
Set-NumericValue $ws.Cells.Item(275, 12) 5                # L275
Set-NumericValue $ws.Cells.Item(275, 13) 3                # M275

# Row 276 (04.12.2020)
$ws.Cells.Item(276, 3).Value = 156                      # C276
Set-NumericValue $ws.Cells.Item(276, 12) 3                # L276

# Row 277 (05.12.2020)
Set-NumericValue $ws.Cells.Item(277, 12) 2                # L277

# Row 278 (06.12.2020)
Set-NumericValue $ws.Cells.Item(278, 13) 1                # M278

# Row 279 (07.12.2020)
$ws.Cells.Item(279, 3).Value = 161                      # C279
Set-NumericValue $ws.Cells.Item(279, 12) 2                # L279
Set-NumericValue $ws.Cells.Item(279, 13) 1                # M279

# Row 280 (08.12.2020) - L/M were previously the text "0"; they become
# real numbers now.
$ws.Cells.Item(280, 3).Value = 93                       # C280
Set-NumericValue $ws.Cells.Item(280, 12) 1                # L280
Set-NumericValue $ws.Cells.Item(280, 13) 0                # M280

# Row 281 (09.12.2020) - previously an all-blank placeholder row; now
# filled in with its own figures.
$ws.Cells.Item(281, 3).Value = 18                       # C281
$ws.Cells.Item(281, 5).Value = 25                       # E281
$ws.Cells.Item(281, 6).Value = 19                       # F281
$ws.Cells.Item(281, 7).Value = 129                      # G281
Set-NumericValue $ws.Cells.Item(281, 12) 2                # L281
Set-NumericValue $ws.Cells.Item(281, 13) 0                # M281

# Row 282 (10.12.2020) remains blank on the input columns; its formulas
# (B/H/J/K) pick up row 281's new totals automatically on recalculation.
